$d = $word.ActiveDocument

# The three trailing paragraphs:
#   # Used by Lightspeed AppBuilder
#   pipeline_type: ls-pipeline-factory
#   one_click_deploy: true
# get replaced by a 15-line Tekton PipelineRun YAML snippet.
# We locate the first of the three paragraphs by its distinctive text,
# then build a Range spanning all three original paragraphs and replace
# their contents (as raw WordprocessingML) in one shot so the run/proofErr
# layout matches exactly.

$w_ns = "http://schemas.openxmlformats.org/wordprocessingml/2006/main"

$startPara = $null
$paras = $d.Paragraphs
for ($i = 1; $i -le $paras.Count; $i++) {
    $txt = $paras.Item($i).Range.Text
    if ($txt -like "*Used by Lightspeed*") {
        $startPara = $i
        break
    }
}

$rng = $d.Range($paras.Item($startPara).Range.Start, $paras.Item($startPara + 2).Range.End)

$frag = @"
<w:p xmlns:w="$w_ns">
  <w:proofErr w:type="spellStart"/>
  <w:r>
    <w:lastRenderedPageBreak/>
    <w:t>apiVersion</w:t>
  </w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r>
    <w:t xml:space="preserve">: </w:t>
  </w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r>
    <w:t>tekton.dev</w:t>
  </w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r>
    <w:t>/v1beta1</w:t>
  </w:r>
</w:p>
<w:p xmlns:w="$w_ns">
  <w:r>
    <w:t xml:space="preserve">kind: </w:t>
  </w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r>
    <w:t>PipelineRun</w:t>
  </w:r>
  <w:proofErr w:type="spellEnd"/>
</w:p>
<w:p xmlns:w="$w_ns">
  <w:r>
    <w:t>metadata:</w:t>
  </w:r>
</w:p>
<w:p xmlns:w="$w_ns">
  <w:r>
    <w:t xml:space="preserve">  </w:t>
  </w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r>
    <w:t>generateName</w:t>
  </w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r>
    <w:t>: ci-sandbox-run-</w:t>
  </w:r>
</w:p>
<w:p xmlns:w="$w_ns">
  <w:r>
    <w:t>spec:</w:t>
  </w:r>
</w:p>
<w:p xmlns:w="$w_ns">
  <w:r>
    <w:t xml:space="preserve">  </w:t>
  </w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r>
    <w:t>pipelineRef</w:t>
  </w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r>
    <w:t>:</w:t>
  </w:r>
</w:p>
<w:p xmlns:w="$w_ns">
  <w:r>
    <w:t xml:space="preserve">    name: ci-sandbox-pipeline</w:t>
  </w:r>
</w:p>
<w:p xmlns:w="$w_ns">
  <w:r>
    <w:t xml:space="preserve">  workspaces:</w:t>
  </w:r>
</w:p>
<w:p xmlns:w="$w_ns">
  <w:r>
    <w:t xml:space="preserve">    - name: shared-workspace</w:t>
  </w:r>
</w:p>
<w:p xmlns:w="$w_ns">
  <w:r>
    <w:t xml:space="preserve">      </w:t>
  </w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r>
    <w:t>volumeClaimTemplate</w:t>
  </w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r>
    <w:t>:</w:t>
  </w:r>
</w:p>
<w:p xmlns:w="$w_ns">
  <w:r>
    <w:t xml:space="preserve">        spec:</w:t>
  </w:r>
</w:p>
<w:p xmlns:w="$w_ns">
  <w:r>
    <w:t xml:space="preserve">          </w:t>
  </w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r>
    <w:t>accessModes</w:t>
  </w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r>
    <w:t>: ["</w:t>
  </w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r>
    <w:t>ReadWriteOnce</w:t>
  </w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r>
    <w:t>"]</w:t>
  </w:r>
</w:p>
<w:p xmlns:w="$w_ns">
  <w:r>
    <w:t xml:space="preserve">          resources:</w:t>
  </w:r>
</w:p>
<w:p xmlns:w="$w_ns">
  <w:r>
    <w:t xml:space="preserve">            requests:</w:t>
  </w:r>
</w:p>
<w:p xmlns:w="$w_ns">
  <w:r>
    <w:t xml:space="preserve">              storage: 1Gi</w:t>
  </w:r>
</w:p>
"@

$rng.InsertXML($frag)
